$d = $word.ActiveDocument

# Locate the target paragraph ("Lens:" answer that begins with
# "Bien sur qu'Il est possible de gagner sa vie en etant developpeur ...")
# by searching for its distinctive opening text, then expanding the
# found range to the whole paragraph. This avoids depending on a
# hard-coded paragraph index.
$rng = $d.Content
$found = $rng.Find.Execute("Bien sûr qu’Il est possible", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Target paragraph not found"
}
$rng.Expand(4) | Out-Null

# Rebuild the paragraph's four runs with the same text/run-level
# metadata as before, but without the explicit w:sz / w:szCs (12pt)
# overrides on the run properties - the size reverts to the document
# default. InsertXML lets us specify the exact resulting OOXML for the
# run properties (dropping <w:sz>/<w:szCs>) while keeping every other
# paragraph/run attribute (paraId, rsids, ...) untouched.
$xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4505A053" w14:textId="77777777" w:rsidR="003419C0" w:rsidRDefault="00434803" w:rsidP="00805410"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00287ADA"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t>Bien s&#251;r qu&#8217;</w:t></w:r><w:r w:rsidR="003A3B9B" w:rsidRPr="00287ADA"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t xml:space="preserve">Il est possible de gagner sa vie en </w:t></w:r><w:r w:rsidR="008D5978" w:rsidRPr="00287ADA"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t>&#233;tant</w:t></w:r><w:r w:rsidR="003A3B9B" w:rsidRPr="00287ADA"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t xml:space="preserve"> d&#233;veloppeur et cela est parmi les meilleurs carri&#232;res pour un informaticien. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xml)
Write-Output "Updated font-size overrides on the 'Bien sur qu'...' paragraph"
